$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I holds free-text counts (e.g. "10") that must remain text, not numbers.
$ws.Range("I59:I77").NumberFormat = "@"

# Row 59
$ws.Range("A59").Value = 112074385
$ws.Range("B59").Value = 78699
$ws.Range("D59").Value = "NT"
$ws.Range("E59").Value = 6458
$ws.Range("F59").Value = "Lunglav"
$ws.Range("G59").Value = "Lobaria pulmonaria"
$ws.Range("H59").Value = "(L.) Hoffm."
$ws.Range("R59").Value = 7054776
$ws.Range("Z59").Value = "18:46"
$ws.Range("AB59").Value = "18:46"

# Row 60
$ws.Range("A60").Value = 112074850
$ws.Range("B60").Value = 96720
$ws.Range("I60").Value = "10"
$ws.Range("Q60").Value = 570507
$ws.Range("R60").Value = 7054761
$ws.Range("Z60").Value = "19:11"
$ws.Range("AB60").Value = "19:11"

# Row 61
$ws.Range("A61").Value = 112074829
$ws.Range("B61").Value = 96720
$ws.Range("I61").Value = "10"
$ws.Range("Q61").Value = 570501
$ws.Range("R61").Value = 7054758
$ws.Range("Z61").Value = "19:10"
$ws.Range("AB61").Value = "19:10"

# Row 62
$ws.Range("A62").Value = 112074371
$ws.Range("B62").Value = 96720
$ws.Range("I62").Value = "100"
$ws.Range("Q62").Value = 570623
$ws.Range("R62").Value = 7054780
$ws.Range("Z62").Value = "18:45"
$ws.Range("AB62").Value = "18:45"

# Row 63
$ws.Range("A63").Value = 112073564
$ws.Range("B63").Value = 96720
$ws.Range("I63").Value = "100"
$ws.Range("Q63").Value = 570581
$ws.Range("R63").Value = 7054735
$ws.Range("Z63").Value = "17:57"
$ws.Range("AB63").Value = "17:57"

# Row 64
$ws.Range("A64").Value = 112074184
$ws.Range("B64").Value = 96720
$ws.Range("I64").Value = "70"
$ws.Range("Q64").Value = 570578
$ws.Range("R64").Value = 7054744
$ws.Range("Z64").Value = "18:35"
$ws.Range("AB64").Value = "18:35"

# Row 65
$ws.Range("A65").Value = 112074315
$ws.Range("B65").Value = 96720
$ws.Range("I65").Value = ""
$ws.Range("Q65").Value = 570621
$ws.Range("R65").Value = 7054778
$ws.Range("Z65").Value = "18:44"
$ws.Range("AB65").Value = "18:44"

# Row 66
$ws.Range("A66").Value = 112074284
$ws.Range("B66").Value = 78699
$ws.Range("E66").Value = 6458
$ws.Range("F66").Value = "Lunglav"
$ws.Range("G66").Value = "Lobaria pulmonaria"
$ws.Range("H66").Value = "(L.) Hoffm."
$ws.Range("M66").Value = ""
$ws.Range("Q66").Value = 570619
$ws.Range("R66").Value = 7054772
$ws.Range("Z66").Value = "18:42"
$ws.Range("AB66").Value = "18:42"

# Row 67
$ws.Range("A67").Value = 112072947
$ws.Range("B67").Value = 96720
$ws.Range("I67").Value = ""
$ws.Range("Q67").Value = 570501
$ws.Range("R67").Value = 7054655
$ws.Range("S67").Value = 1
$ws.Range("Z67").Value = "17:31"
$ws.Range("AB67").Value = "17:31"

# Row 68
$ws.Range("A68").Value = 112073383
$ws.Range("B68").Value = 96720
$ws.Range("D68").Value = "VU"
$ws.Range("E68").Value = 220787
$ws.Range("F68").Value = "Knärot"
$ws.Range("G68").Value = "Goodyera repens"
$ws.Range("H68").Value = "(L.) R. Br."
$ws.Range("I68").Value = "200"
$ws.Range("Q68").Value = 570562
$ws.Range("R68").Value = 7054716
$ws.Range("Z68").Value = "17:50"
$ws.Range("AB68").Value = "17:50"

# Row 69
$ws.Range("A69").Value = 112074220
$ws.Range("B69").Value = 96720
$ws.Range("D69").Value = "VU"
$ws.Range("E69").Value = 220787
$ws.Range("F69").Value = "Knärot"
$ws.Range("G69").Value = "Goodyera repens"
$ws.Range("H69").Value = "(L.) R. Br."
$ws.Range("I69").Value = "20"
$ws.Range("Q69").Value = 570573
$ws.Range("R69").Value = 7054742
$ws.Range("Z69").Value = "18:37"
$ws.Range("AB69").Value = "18:37"

# Row 70
$ws.Range("A70").Value = 112072920
$ws.Range("B70").Value = 96720
$ws.Range("I70").Value = "60"
$ws.Range("Q70").Value = 570486
$ws.Range("R70").Value = 7054643
$ws.Range("Z70").Value = "17:29"
$ws.Range("AB70").Value = "17:29"

# Row 71
$ws.Range("A71").Value = 112073748
$ws.Range("B71").Value = 96720
$ws.Range("Q71").Value = 570532
$ws.Range("R71").Value = 7054761
$ws.Range("Z71").Value = "18:09"
$ws.Range("AB71").Value = "18:09"

# Row 72
$ws.Range("A72").Value = 112073661
$ws.Range("B72").Value = 56430
$ws.Range("D72").Value = "NT"
$ws.Range("E72").Value = 100109
$ws.Range("F72").Value = "Tretåig hackspett"
$ws.Range("G72").Value = "Picoides tridactylus"
$ws.Range("H72").Value = "(Linnaeus, 1758)"
$ws.Range("I72").Value = ""
$ws.Range("M72").Value = "äldre spår"
$ws.Range("Q72").Value = 570523
$ws.Range("R72").Value = 7054775
$ws.Range("S72").Value = 25
$ws.Range("Z72").Value = "18:03"
$ws.Range("AB72").Value = "18:03"

# Row 73
$ws.Range("A73").Value = 112073706
$ws.Range("B73").Value = 96720
$ws.Range("D73").Value = "VU"
$ws.Range("E73").Value = 220787
$ws.Range("F73").Value = "Knärot"
$ws.Range("G73").Value = "Goodyera repens"
$ws.Range("H73").Value = "(L.) R. Br."
$ws.Range("I73").Value = "10"
$ws.Range("M73").Value = ""
$ws.Range("Q73").Value = 570517
$ws.Range("R73").Value = 7054754
$ws.Range("Z73").Value = "18:06"
$ws.Range("AB73").Value = "18:06"

# Row 74
$ws.Range("A74").Value = 112073635
$ws.Range("B74").Value = 96720
$ws.Range("I74").Value = "30"
$ws.Range("Q74").Value = 570513
$ws.Range("R74").Value = 7054747
$ws.Range("Z74").Value = "18:01"
$ws.Range("AB74").Value = "18:01"

# Row 75
$ws.Range("A75").Value = 112073927
$ws.Range("B75").Value = 56430
$ws.Range("D75").Value = "NT"
$ws.Range("E75").Value = 100109
$ws.Range("F75").Value = "Tretåig hackspett"
$ws.Range("G75").Value = "Picoides tridactylus"
$ws.Range("H75").Value = "(Linnaeus, 1758)"
$ws.Range("I75").Value = ""
$ws.Range("M75").Value = "äldre spår"
$ws.Range("Q75").Value = 570548
$ws.Range("R75").Value = 7054724
$ws.Range("Z75").Value = "18:19"
$ws.Range("AB75").Value = "18:19"

# Row 76
$ws.Range("A76").Value = 112074007
$ws.Range("B76").Value = 96720
$ws.Range("Q76").Value = 570552
$ws.Range("R76").Value = 7054717
$ws.Range("Z76").Value = "18:25"
$ws.Range("AB76").Value = "18:25"

# Row 77
$ws.Range("A77").Value = 112074296
$ws.Range("B77").Value = 96720
$ws.Range("I77").Value = "150"
$ws.Range("Q77").Value = 570620
$ws.Range("R77").Value = 7054773
$ws.Range("Z77").Value = "18:43"
$ws.Range("AB77").Value = "18:43"

